$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new row 6 with a multi-character string value "abc"
$ws.Range("A6").Value = "abc"

# Move the active selection to A7 (matches post-edit sheet view state)
$ws.Range("A7").Select()
